$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are text-formatted in the source data
# (values like "69.334.05" or "0.0000351" must stay literal text, not be
# coerced to numbers -- set NumberFormat to Text before writing the values).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.288.98'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.942.88'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '495.53'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").Value = '147.95'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  -0.99%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("E10").Value = '  +3.84%  '
$ws.Range("D11").Value = '0.0000350'
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").Value = '43.38'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '10.47'
$ws.Range("E13").Value = '  -2.63%  '
$ws.Range("D14").Value = '4.575.62'
$ws.Range("D15").Value = '3.959.62'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '14.23'
$ws.Range("E16").Value = '  -3.48%  '
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  +4.60%  '
$ws.Range("D19").Value = '19.95'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("D20").Value = '69.363.17'
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("D21").Value = '438.12'
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D23").Value = '14.64'
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("D24").Value = '88.89'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '12.03'
$ws.Range("E25").Value = '  +5.80%  '
$ws.Range("E26").Value = '  +4.06%  '
$ws.Range("D27").Value = '11.14'
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").Value = '37.18'
$ws.Range("E28").Value = '  -4.61%  '
$ws.Range("D29").Value = '5.65'
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").Value = '700.09'
$ws.Range("E30").Value = '  -3.05%  '
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").Value = '2.86'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").Value = '0.460'
$ws.Range("E34").Value = '  +14.26%  '

# Rows 35/36 swapped rank order: PEPE (was row 35) and OKB (was row 36)
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '62.57'
$ws.Range("E35").Value = '  +2.97%  '
$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0892'
$ws.Range("E36").Value = '  -2.46%  '

$ws.Range("D37").Value = '6.08'
$ws.Range("E37").Value = '  -1.67%  '
$ws.Range("D38").Value = '41.08'
$ws.Range("E38").Value = '  -2.75%  '
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").Value = '0.997'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '0.0488'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").Value = '2.91'
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").Value = '3.38'
$ws.Range("E47").Value = '  +6.24%  '
$ws.Range("D48").Value = '2.99'
$ws.Range("E48").Value = '  +4.58%  '
$ws.Range("D49").Value = '3.39'
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("D50").Value = '0.0₆0350'
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("E51").Value = '  -2.92%  '
